# Update cryptocurrency price/volume data on Sheet1 (A1:E51 table).
# Also swaps the Stacks/OKB rows (43 and 44).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.704.68"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "2.404.76"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'565.72"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "'137.84"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "2.382.61"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "'5.02"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "'0.335"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").Value = "'25.72"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'0.0000169"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "60.725.67"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "2.399.72"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "'7.76"
$ws.Range("E19").Value = "  +8.20%  "
$ws.Range("D20").Value = "'10.53"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "'320.51"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'4.00"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("D23").Value = "'6.07"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'1.80"
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").Value = "'64.69"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").Value = "'571.60"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("D28").Value = "'8.16"
$ws.Range("E28").Value = "  -10.07%  "
$ws.Range("D29").Value = "2.536.33"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").Value = "0.0₃0909"
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("D31").Value = "'7.76"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "'1.33"
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  -4.77%  "
$ws.Range("D34").Value = "'0.130"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'151.40"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "'1.38"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").Value = "'4.54"
$ws.Range("E38").Value = "  -6.97%  "
$ws.Range("D39").Value = "'0.363"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").Value = "'18.05"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").Value = "'5.04"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'41.07"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.64"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("D45").Value = "'2.24"
$ws.Range("E45").Value = "  -10.41%  "
$ws.Range("D46").Value = "'140.91"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "0.0₆0259"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").Value = "'3.46"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("D49").Value = "'0.581"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "'0.0496"
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("D51").Value = "'19.09"
$ws.Range("E51").Value = "  -3.50%  "
